# Refresh the crypto price/volume table (columns D "Price" and E
# "Volume(1h)") on the active worksheet with newly scraped figures,
# matching the upstream GitHub Actions data-refresh commit.
#
# Column D holds price text that sometimes looks numeric (e.g. "211.56",
# "1.00"); NumberFormat is forced to Text ("@") first for those cells so
# Excel keeps them as literal strings (preserving trailing zeros like
# "1.00") instead of silently converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.873.46'
$ws.Range("E2").Value = '  -0.27%  '
$ws.Range("D3").Value = '1.628.86'
$ws.Range("E3").Value = '  -0.09%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.56'
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.516'
$ws.Range("E6").Value = '  -1.15%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.25'
$ws.Range("E8").Value = '  -0.62%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.257'
$ws.Range("E9").Value = '  -0.26%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0607'
$ws.Range("E10").Value = '  -1.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0879'
$ws.Range("E11").Value = '  -0.12%  '
$ws.Range("D12").Value = '1.860.71'
$ws.Range("E12").Value = '  -0.02%  '
$ws.Range("D13").Value = '1.636.54'
$ws.Range("E13").Value = '  +0.95%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.01'
$ws.Range("E14").Value = '  -0.85%  '
$ws.Range("E15").Value = '  -1.09%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.89'
$ws.Range("E16").Value = '  -1.07%  '
$ws.Range("D17").Value = '27.900.63'
$ws.Range("E17").Value = '  -0.13%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '227.74'
$ws.Range("E18").Value = '  -1.29%  '
$ws.Range("E19").Value = '  -0.28%  '
$ws.Range("D20").Value = '0.0₃0718'
$ws.Range("E20").Value = '  -0.92%  '
$ws.Range("E21").Value = '  -0.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.33'
$ws.Range("E22").Value = '  -0.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.95'
$ws.Range("E23").Value = '  -4.11%  '
$ws.Range("E24").Value = '  +0.67%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '155.31'
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("E26").Value = '  -0.43%  '
$ws.Range("E27").Value = '  -0.45%  '
$ws.Range("E28").Value = '  -0.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.40'
$ws.Range("E29").Value = '  -1.13%  '
$ws.Range("E30").Value = '  -0.41%  '
$ws.Range("E31").Value = '  -0.33%  '
$ws.Range("E32").Value = '  +0.12%  '
$ws.Range("D33").Value = '1.413.16'
$ws.Range("E33").Value = '  +1.04%  '
$ws.Range("E34").Value = '  +0.97%  '
$ws.Range("E35").Value = '  +2.79%  '
$ws.Range("E36").Value = '  -3.93%  '
$ws.Range("E37").Value = '  -1.34%  '
$ws.Range("E38").Value = '  -0.98%  '
$ws.Range("E39").Value = '  -0.55%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.851'
$ws.Range("E40").Value = '  -1.64%  '
$ws.Range("E41").Value = '  -1.79%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '65.74'
$ws.Range("E42").Value = '  -1.13%  '
$ws.Range("E43").Value = '  -0.64%  '
$ws.Range("E44").Value = '  -0.77%  '
$ws.Range("D45").Value = '1.770.07'
$ws.Range("E45").Value = '  -0.10%  '
$ws.Range("E46").Value = '  -3.66%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '88.67'
$ws.Range("E47").Value = '  +0.60%  '
$ws.Range("E48").Value = '  +1.19%  '
$ws.Range("E49").Value = '  -0.36%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.59'
$ws.Range("E50").Value = '  +0.91%  '
$ws.Range("E51").Value = '  -0.04%  '
